$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it occurs ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        while ($true) {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
            if ($found -eq $null -or $found.Address() -eq $firstAddress) {
                break
            }
        }
    }
}

# --- Narrow the "zh-cn"/"de-de" status columns on Overview, and "Status" column on the
#     per-locale sheets from ~17.22 chars to ~13.41 chars ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
